{"js": "// The document contains a run of text that spells out \"<id>p136v_1</id>\"\n// across three separate runs: \"<id>\", \"p136v_1\", \"</id>\". The edit merges\n// them into a single run containing the full string \"<id>p136v_1</id>\"\n// (taking on the formatting of the first of the three runs).\n//\n// We find that exact contiguous text and rewrite it in place, which causes\n// the engine to collapse it back down to one run using the leading run's\n// formatting - matching the target OOXML.\nconst target = \"<id>p136v_1</id>\";\n\nconst results = context.document.body.search(target, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const range = results.items[0];\n  range.insertText(target, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# The document spells \"<id>p136v_1</id>\" across three separate runs:\n# \"<id>\", \"p136v_1\", \"</id>\" (the middle run carries different formatting\n# than the two tag runs around it). The edit merges them into a single run\n# containing the full string \"<id>p136v_1</id>\", using the formatting of the\n# first (\"<id>\") run.\n#\n# Using Find & Replace with the exact same text as both the search target\n# and the replacement causes Word to collapse the matched span back down to\n# a single run (seeded from the first run's formatting), which reproduces\n# the target edit.\n\n$d = $word.ActiveDocument\n$target = \"<id>p136v_1</id>\"\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Text = $target\n$range.Find.Replacement.ClearFormatting()\n$range.Find.Replacement.Text = $target\n$range.Find.Forward = $true\n$range.Find.Wrap = 1          # wdFindContinue\n$range.Find.Format = $false\n$range.Find.MatchCase = $true\n$range.Find.MatchWholeWord = $false\n$range.Find.MatchWildcards = $false\n\n$found = $range.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $target, 1)\nWrite-Output \"found: $found\"\n"}
